# Insert a new data row at row 141 (pushes existing rows 141-209 down to 142-210),
# duplicate the now-shifted row 142 into the new row 141, then overwrite the
# columns that hold the new record's distinct values (Fecha, Volumen,
# Precio minimo/maximo/promedio ponderado, Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(141).Insert()

$ws.Range("A142:R142").Copy()
$ws.Range("A141").PasteSpecial()

$ws.Range("D141").Value = 44460
$ws.Range("J141").Value = 35
$ws.Range("K141").Value = 33000
$ws.Range("L141").Value = 33000
$ws.Range("M141").Value = 33000
$ws.Range("P141").Value = 1320
